$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Remove the naive-forecaster component cells that were buggy ---
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("C4").ClearContents()

# --- Recalculated values after fixing the naive component forecaster bug ---
$ws.Range("E4").Value = -0.6367039903685923
$ws.Range("C5").Value = -4.857241224140941
$ws.Range("C6").Value = -3.956152295564885
$ws.Range("C7").Value = 0.3625742673738719
$ws.Range("C8").Value = 1.234995474941436
$ws.Range("E8").Value = -0.209816187756795
$ws.Range("C10").Value = 0.8993608108207818
$ws.Range("C11").Value = 0.6836026627130787
$ws.Range("C13").Value = 0.2336505480021511
$ws.Range("C14").Value = 0.02019328874802717
$ws.Range("E15").Value = -0.03923323971221082
$ws.Range("C17").Value = -0.1588690085688071
$ws.Range("E17").Value = -0.4617525814883061
$ws.Range("E18").Value = -0.001769149545449711
$ws.Range("C19").Value = -0.5438176183081955
$ws.Range("E19").Value = 0.01247916696662799
$ws.Range("C21").Value = -0.006876704825686808
$ws.Range("C22").Value = 0.0720185131838802
$ws.Range("C23").Value = 0.6840863075407544
$ws.Range("E26").Value = -0.950584780912811
$ws.Range("C27").Value = -0.8803054679952238
$ws.Range("E27").Value = -0.2797061371760057
$ws.Range("C28").Value = -0.63478973259814
$ws.Range("E28").Value = 0.01241557525979431
$ws.Range("C29").Value = -0.8089889044073373
$ws.Range("E29").Value = -0.2234428210500905
$ws.Range("C30").Value = -0.8017595264762423
$ws.Range("E30").Value = 0.0476740348578808
$ws.Range("C31").Value = -0.3388987799285315
$ws.Range("E31").Value = -0.1599040255974238
$ws.Range("C32").Value = -0.7158018152081613
$ws.Range("E32").Value = -0.7615805088034833
$ws.Range("C36").Value = 1.576357831383679
$ws.Range("C37").Value = 0.9319769131821865
$ws.Range("E37").Value = -0.7143539453189907
$ws.Range("C38").Value = 0.9704846793491706
$ws.Range("E38").Value = -0.8754609427830351
$ws.Range("C39").Value = -0.05255054327948372
$ws.Range("C40").Value = -1.488707312182613
$ws.Range("E40").Value = -1.062239424572287
$ws.Range("C41").Value = -1.002933118524785
$ws.Range("E41").Value = -0.8142357075366813
$ws.Range("E45").Value = -0.5312671580391171
$ws.Range("C46").Value = 0.3928252664241683
$ws.Range("E47").Value = 0.1200540108007964
$ws.Range("C48").Value = 0.8934739937295433
$ws.Range("E48").Value = 0.4129745242491101
$ws.Range("C49").Value = 0.7146010879610643
$ws.Range("C50").Value = 0.3224026462283369
$ws.Range("E50").Value = -0.9749878381046684
$ws.Range("C51").Value = -2.355806475145761
$ws.Range("C52").Value = -1.372720900450863
